# Problem Set 8-wk08.docx edit
#
# 1) Merge the three runs in the "HINT: For problems 5-7 ..." paragraph
#    into a single run with the combined text (the split "a" / "nd their
#    solutions..." runs go away, the leading run absorbs the whole
#    sentence).
# 2) Flip the "Normal" style's overflow-punctuation paragraph setting
#    from On to Off (w:overflowPunct true -> false).

$d = $word.ActiveDocument

$old = ": For problems 5 -7, look at the quick check problems from page 858 and their solutions on page 880 as a template."
$new = ": For problems 5 -7, look at the quick check problems from page 858 and their solutions on page 880 as a template."

# Word's Find/Replace re-writes the matched range as a single run, which
# collapses the original 3-run split ("...page 858 " + "a" + "nd their...")
# into one run carrying the full sentence.
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

# The built-in "Normal" paragraph style no longer forces hanging/overflow
# punctuation.
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.HangingPunctuation = $false
